$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new error code numbers to existing rows 5 and 6
$ws.Range("A5").Value = 103
$ws.Range("A6").Value = 104

# Add new rows 7-10 with new error codes
$ws.Range("A7").Value = 105
$ws.Range("B7").Value = "VERSION_INVALID"
$ws.Range("C7").Value = "V zadanom zázname neexistuje verzia"

$ws.Range("A8").Value = 106
$ws.Range("B8").Value = "VERSION_UNAVAILABLE"
$ws.Range("C8").Value = "Zadaná verzia nespĺňa kritéria pre SW ver. 2G ani 3G"

$ws.Range("A9").Value = 107
$ws.Range("B9").Value = "NO_PROCESSED_RECORDS"
$ws.Range("C9").Value = "Žiaden zo záznamov sa nepodarilo spracovať"
$ws.Range("D9").Value = "Kontrola formátu záznamu príp. Jeho úprava"

$ws.Range("A10").Value = 108
$ws.Range("B10").Value = "MISSING_SAFE_BITES"
$ws.Range("C10").Value = "V zázname neboli nájdené safe bytes"
$ws.Range("D10").Value = "Doplniť chýbajúce safe bytes"

# Resize the table to include the new rows
$tbl = $ws.ListObjects.Item("Table1")
$tbl.Resize($ws.Range("A1:D10"))

# Set column B width to match bestFit width behavior (Excel-computed best fit)
$ws.Columns.Item(2).ColumnWidth = 24.5703125

# Update selection to match final state
$ws.Range("C10").Select()
